# metodosnum2020.xlsx — "update excel - hmk 1,2"
#
# Marks homeworks 1 & 2 ("tarea1 (scripts)" / "tarea2 (finacci, factorial)")
# as turned in (F/G columns) and graded (H column) for every student whose
# row previously had 0s, gives four students (rows 4, 8, 15, 18, 21, 23)
# their first score entries, fixes a couple of half-credit (0.5) grades up
# to full credit, flips I33 back to 0, and leaves a grading note in J6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => { Column => new value } for every touched grading cell (F, G, H, I)
$updates = @(
    @{ Row = 3;  Cells = @{ "F" = 1; "G" = 1 } }
    @{ Row = 4;  Cells = @{ "F" = 0; "G" = 0; "H" = 0; "I" = 0 } }
    @{ Row = 5;  Cells = @{ "F" = 1; "G" = 1; "H" = 1 } }
    @{ Row = 6;  Cells = @{ "F" = 1; "G" = 1; "H" = 1 } }
    @{ Row = 7;  Cells = @{ "F" = 1; "G" = 1; "H" = 1 } }
    @{ Row = 8;  Cells = @{ "F" = 0; "G" = 0; "H" = 0; "I" = 0 } }
    @{ Row = 9;  Cells = @{ "F" = 1; "G" = 1; "H" = 1 } }
    @{ Row = 10; Cells = @{ "F" = 1; "G" = 1; "H" = 1 } }
    @{ Row = 11; Cells = @{ "F" = 1; "G" = 1; "H" = 1 } }
    @{ Row = 12; Cells = @{ "F" = 1; "G" = 1; "H" = 1 } }
    @{ Row = 13; Cells = @{ "F" = 1; "G" = 1; "H" = 1 } }
    @{ Row = 14; Cells = @{ "F" = 1; "G" = 1; "H" = 1 } }
    @{ Row = 15; Cells = @{ "F" = 1; "G" = 1; "H" = 1; "I" = 0 } }
    @{ Row = 17; Cells = @{ "F" = 1; "G" = 1; "H" = 1 } }
    @{ Row = 18; Cells = @{ "F" = 1; "G" = 1; "H" = 1; "I" = 0 } }
    @{ Row = 19; Cells = @{ "H" = 1 } }
    @{ Row = 20; Cells = @{ "F" = 1; "G" = 1; "H" = 1 } }
    @{ Row = 21; Cells = @{ "F" = 1; "G" = 1; "H" = 1; "I" = 0 } }
    @{ Row = 22; Cells = @{ "F" = 1; "G" = 1; "H" = 1 } }
    @{ Row = 23; Cells = @{ "F" = 1; "G" = 1; "H" = 1; "I" = 0 } }
    @{ Row = 24; Cells = @{ "F" = 1; "G" = 1; "H" = 1 } }
    @{ Row = 25; Cells = @{ "F" = 1; "G" = 1; "H" = 1 } }
    @{ Row = 26; Cells = @{ "F" = 1; "G" = 1; "H" = 1 } }
    @{ Row = 27; Cells = @{ "F" = 1; "G" = 1; "H" = 1 } }
    @{ Row = 28; Cells = @{ "F" = 1; "G" = 1; "H" = 1 } }
    @{ Row = 29; Cells = @{ "F" = 1; "G" = 1; "H" = 1 } }
    @{ Row = 30; Cells = @{ "F" = 1; "G" = 1; "H" = 1 } }
    @{ Row = 32; Cells = @{ "F" = 1; "G" = 1; "H" = 1 } }
    @{ Row = 33; Cells = @{ "F" = 1; "G" = 1; "H" = 1; "I" = 0 } }
    @{ Row = 34; Cells = @{ "F" = 1; "G" = 1 } }
    @{ Row = 35; Cells = @{ "F" = 1; "G" = 1; "H" = 1 } }
    @{ Row = 37; Cells = @{ "F" = 1; "G" = 1; "H" = 1 } }
    @{ Row = 38; Cells = @{ "F" = 1; "G" = 1; "H" = 1 } }
)

foreach ($u in $updates) {
    $row = $u.Row
    foreach ($col in $u.Cells.Keys) {
        $ws.Range("$col$row").Value = $u.Cells[$col]
    }
}

# Grading note added for row 6 (new column J)
$ws.Range("J6").Value = "se deben subir los .m"

# Move the saved selection/scroll position to E8 (also drops the stale
# topLeftCell="A10" the sheet view previously had).
$ws.Range("E8").Select()
